$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D/E are plain text in the source data. Some of the updated Price
# values look like numbers (e.g. "580.40", "12.00"); force those specific
# cells to Text format first so Excel keeps the exact string (trailing
# zeros, etc.) instead of silently converting them to a Number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "64.242.73"
$ws.Range("E2").Value = "  +5.31%  "
$ws.Range("D3").Value = "2.763.20"
$ws.Range("E3").Value = "  +3.38%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "580.40"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "155.04"
$ws.Range("E6").Value = "  +6.80%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("D9").Value = "2.759.54"
$ws.Range("E9").Value = "  +3.29%  "
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("E11").Value = "  +4.54%  "
$ws.Range("D12").Value = "0.389"
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("E13").Value = "  +3.53%  "
$ws.Range("D14").Value = "3.250.39"
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("D15").Value = "26.52"
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("D16").Value = "64.127.68"
$ws.Range("E16").Value = "  +5.15%  "
$ws.Range("D17").Value = "0.0000153"
$ws.Range("E17").Value = "  +6.32%  "
$ws.Range("D18").Value = "2.757.90"
$ws.Range("E18").Value = "  +3.63%  "
$ws.Range("D19").Value = "12.00"
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("E20").Value = "  +2.46%  "
$ws.Range("D21").Value = "361.74"
$ws.Range("E21").Value = "  +2.95%  "
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "0.533"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").Value = "66.52"
$ws.Range("E25").Value = "  +4.01%  "
$ws.Range("E26").Value = "  +4.96%  "
$ws.Range("D27").Value = "8.53"
$ws.Range("E27").Value = "  +4.59%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").Value = "0.0₃0914"
$ws.Range("E29").Value = "  +11.85%  "
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("E31").Value = "  +2.91%  "
$ws.Range("E32").Value = "  +16.48%  "
$ws.Range("D33").Value = "172.13"
$ws.Range("E33").Value = "  +4.00%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "20.46"
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("D36").Value = "4.84"
$ws.Range("E36").Value = "  +7.85%  "
$ws.Range("E37").Value = "  +8.17%  "
$ws.Range("D38").Value = "1.82"
$ws.Range("E38").Value = "  +9.42%  "
$ws.Range("D39").Value = "1.01"
$ws.Range("E39").Value = "  +14.70%  "
$ws.Range("D40").Value = "347.66"
$ws.Range("E40").Value = "  +4.05%  "
$ws.Range("D41").Value = "4.25"
$ws.Range("E41").Value = "  +5.36%  "
$ws.Range("D42").Value = "39.20"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("D43").Value = "5.73"
$ws.Range("E43").Value = "  +9.57%  "
$ws.Range("D44").Value = "21.85"
$ws.Range("E44").Value = "  +6.68%  "
$ws.Range("D45").Value = "21.73"
$ws.Range("E45").Value = "  +5.92%  "
$ws.Range("D46").Value = "0.0591"
$ws.Range("E46").Value = "  +4.72%  "
$ws.Range("D47").Value = "0.646"
$ws.Range("E47").Value = "  +4.93%  "
$ws.Range("D48").Value = "137.35"
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("D49").Value = "0.0256"
$ws.Range("E49").Value = "  +2.58%  "
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("E51").Value = "  +0.07%  "
